# v1.2 - reopen task again to modify something
# (ask to add more details to registration form wireframe)

$wb       = $excel.ActiveWorkbook
$wsLogin  = $wb.Worksheets.Item("LH_Review_WF_LOGIN")
$wsHist   = $wb.Worksheets.Item("VERSION-HISTORY")

# --- LH_Review_WF_LOGIN: reopen all three review items ---
# "Reviewer verification" column (I) goes back from "closed" to "open"
$wsLogin.Range("I2").Value = "open"
$wsLogin.Range("I3").Value = "open"
$wsLogin.Range("I4").Value = "open"

# --- VERSION-HISTORY: the v1.1 "closed" entry is undone/cleared ---
$wsHist.Range("A3:D3").ClearContents()
$wsHist.Range("A3").Interior.Color = 16382198
$wsHist.Range("D3").Interior.Color = 16382198

# --- Restore selections / active sheet bookkeeping ---
$wsHist.Range("C19").Select()
$wsLogin.Activate()
$wsLogin.Range("I11").Select()
